# Employee records now replace the old traffic-incident report schema.
# Strategy:
#   1. Drop the two trailing "responders.personnel.*" columns (R:S) and the
#      two now-empty trailing rows (8:9) so the sheet shrinks to A1:Q7.
#   2. Force the date-like text columns (C, I, O, Q) to be stored as literal
#      text instead of being auto-parsed into date serial numbers.
#   3. Overwrite every header and every data cell that actually changed.
#   4. performanceReviews.metrics.scores (column K) holds real numbers now,
#      so those six cells are written as numeric values, not strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shrink the sheet to the new used range -----------------------------
$ws.Columns("R:S").Delete()
$ws.Rows("8:9").Delete()

# --- 2. Keep YYYY-MM-DD-looking strings as text, not dates ------------------
$ws.Range("C2:C7").NumberFormat = "@"
$ws.Range("I2:I7").NumberFormat = "@"
$ws.Range("O2:O7").NumberFormat = "@"
$ws.Range("Q2:Q7").NumberFormat = "@"

# --- 3. Header row -----------------------------------------------------------
$ws.Range("A1").Value = "employeeId"
$ws.Range("B1").Value = "firstName"
$ws.Range("C1").Value = "personalInfo.dateOfBirth"
$ws.Range("D1").Value = "personalInfo.phoneNumber"
$ws.Range("E1").Value = "personalInfo.emergencyContact"
$ws.Range("F1").Value = "lastName"
$ws.Range("G1").Value = "email"
$ws.Range("H1").Value = "performanceReviews.reviewId"
$ws.Range("I1").Value = "performanceReviews.reviewDate"
$ws.Range("J1").Value = "performanceReviews.metrics.metricName"
$ws.Range("K1").Value = "performanceReviews.metrics.scores"
$ws.Range("L1").Value = "skills.skillName"
$ws.Range("M1").Value = "skills.proficiencyLevel"
$ws.Range("N1").Value = "trainingHistory.courseName"
$ws.Range("O1").Value = "trainingHistory.completionDate"
$ws.Range("P1").Value = "trainingHistory.certifications.certificationName"
$ws.Range("Q1").Value = "trainingHistory.certifications.issueDate"

# --- 4. Row 2 (employee EMP-2024-001, Sarah Johnson) ------------------------
$ws.Range("A2").Value = "EMP-2024-001"
$ws.Range("B2").Value = "Sarah"
$ws.Range("C2").Value = "1990-03-15"
$ws.Range("D2").Value = "+1-555-0123"
$ws.Range("E2").Value = "John Mitchell (Spouse)"
$ws.Range("F2").Value = "Johnson"
$ws.Range("G2").Value = "sarah.johnson@company.com"
$ws.Range("H2").Value = "REV-2024-Q1"
$ws.Range("I2").Value = "2024-03-15"
$ws.Range("J2").Value = "Communication 1"
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = "JavaScript"
$ws.Range("M2").Value = "Expert"
$ws.Range("N2").Value = "Advanced Leadership Development"
$ws.Range("O2").Value = "2024-02-28"
$ws.Range("P2").Value = "Certified Scrum Master"
$ws.Range("Q2").Value = "2024-01-15"

# --- 5. Row 3 (continuation: skills / training / certifications) ------------
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "Communication 2"
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = "Project Management"
$ws.Range("M3").Value = "Intermediate"
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = ""
$ws.Range("P3").Value = "AWS Solutions Architect"
$ws.Range("Q3").Value = "2023-11-20"

# --- 6. Row 4 (second performance review) ------------------------------------
$ws.Range("H4").Value = "REV-2024-Q2"
$ws.Range("I4").Value = "2024-03-15"
$ws.Range("J4").Value = "Communication 3"
$ws.Range("K4").Value = 2
$ws.Range("P4").Value = ""
$ws.Range("Q4").Value = ""

# --- 7. Row 5 (extra metric row) ---------------------------------------------
$ws.Range("J5").Value = "Communication 4"
$ws.Range("K5").Value = 3

# --- 8. Row 6 (employee EMP-2024-002, Michael Rodriguez) --------------------
$ws.Range("A6").Value = "EMP-2024-002"
$ws.Range("B6").Value = "Michael"
$ws.Range("C6").Value = "1985-07-22"
$ws.Range("D6").Value = "+1-555-0456"
$ws.Range("E6").Value = "Lisa Chen (Sister)"
$ws.Range("F6").Value = "Rodriguez"
$ws.Range("G6").Value = "michael.rodriguez@company.com"
$ws.Range("H6").Value = "REV-2024-Q1-002"
$ws.Range("I6").Value = "2024-03-20"
$ws.Range("J6").Value = "Technical Skills 1"
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = "Employee Relations"
$ws.Range("M6").Value = "Expert"
$ws.Range("N6").Value = "Diversity and Inclusion Workshop"
$ws.Range("O6").Value = "2024-01-10"
$ws.Range("P6").Value = "PHR Certification"
$ws.Range("Q6").Value = "2023-08-15"

# --- 9. Row 7 (continuation row) ---------------------------------------------
$ws.Range("H7").Value = "REV-2024-Q1-003"
$ws.Range("I7").Value = "2024-03-20"
$ws.Range("J7").Value = "Technical Skills 2"
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = "Data Analysis"
$ws.Range("M7").Value = "Intermediate"
$ws.Range("N7").Value = ""
$ws.Range("O7").Value = ""
